# Week_2.docx edits: correct a handful of wording errors and drop the
# stray "References" section (page-break + heading + editor's question)
# that was left over at the end of the document, leaving only the
# bookmark that used to sit inside that paragraph.

$d = $word.ActiveDocument

# 1) "на поряд простіше" -> "на порядок простіше"
$d.Content.Find.Execute("на поряд простіше", $true, $false, $false, $false, $false, `
    $true, 1, $false, "на порядок простіше", 2) | Out-Null

# 2) "«парувати»" -> "парити"
$quote1 = [char]0x00AB
$quote2 = [char]0x00BB
$d.Content.Find.Execute("$quote1" + "парувати" + "$quote2", $true, $false, $false, $false, $false, `
    $true, 1, $false, "парити", 2) | Out-Null

# 3) "не повинна ставати на значно впливати" -> "не повинна значно впливати"
$d.Content.Find.Execute("не повинна ставати на значно впливати", $true, $false, $false, $false, $false, `
    $true, 1, $false, "не повинна значно впливати", 2) | Out-Null

# 4) Remove the trailing page break + "Список використаної літератури"
#    heading paragraph, and strip the leftover editor's comment
#    paragraph down to just its bookmark, marking that paragraph en-US.
$target = "А список літератури треба до кожного розділу? Як в Тустанівського? Це ж розрахункова, хіба треба до кожного розділу писати окремий список?"

$paras = $d.Paragraphs
$qIndex = $paras.Count
$qPara = $paras.Item($qIndex)
while ($qPara.Range.Text.IndexOf($target) -lt 0 -and $qIndex -gt 1) {
    $qIndex = $qIndex - 1
    $qPara = $paras.Item($qIndex)
}

$headingPara = $paras.Item($qIndex - 1)
$pageBreakPara = $paras.Item($qIndex - 2)

$rRemove = $d.Range($pageBreakPara.Range.Start, $headingPara.Range.End)
$rRemove.Delete()

# The structural edit above invalidates old paragraph/range handles -
# re-fetch the (now-shifted) question paragraph from the collection.
$newCount = $d.Paragraphs.Count
$qPara = $d.Paragraphs.Item($newCount)

$fr = $qPara.Range.Duplicate
$found = $fr.Find.Execute($target, $true)
if ($found) {
    $fr.Delete()
}

$d.Paragraphs.Item($d.Paragraphs.Count).Range.LanguageID = "en-US"
